$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 3785.8462
$ws.Range("I34").Value = 3618
$ws.Range("K34").Value = 3618
$ws.Range("M34").Value = -3415
$ws.Range("H36").Value = 3785.8462
$ws.Range("I36").Value = 3618
$ws.Range("K36").Value = 3618
$ws.Range("M36").Value = -2903
$ws.Range("H70").Value = 2580.5
$ws.Range("I70").Value = 3240.75
$ws.Range("J70").Value = 2250.375
$ws.Range("K70").Value = 9722.25
$ws.Range("L70").Value = 6751.125
$ws.Range("M70").Value = -9452.25
$ws.Range("N70").Value = -7291.125
$ws.Range("H73").Value = 2580.5
$ws.Range("I73").Value = 3240.75
$ws.Range("J73").Value = 2250.375
$ws.Range("K73").Value = 9722.25
$ws.Range("L73").Value = 6751.125
$ws.Range("M73").Value = -8786.25
$ws.Range("N73").Value = -8623.125
$ws.Range("H94").Value = 831.7
$ws.Range("I94").Value = 831.7
$ws.Range("K94").Value = 831.7
$ws.Range("M94").Value = -380.7
$ws.Range("H98").Value = 1693.9231
$ws.Range("J98").Value = 5077.8
$ws.Range("L98").Value = 5077.8
$ws.Range("N98").Value = -8073.8
$ws.Range("H116").Value = 7412.1724
$ws.Range("I116").Value = 9444.277
$ws.Range("J116").Value = 4086.9092
$ws.Range("K116").Value = 9444.277
$ws.Range("L116").Value = 4086.9092
$ws.Range("M116").Value = -6002.277
$ws.Range("N116").Value = -10970.9092
$ws.Range("H122").Value = 1693.9231
$ws.Range("J122").Value = 5077.8
$ws.Range("L122").Value = 15233.4
$ws.Range("N122").Value = -20133.4
$ws.Range("H123").Value = 100000
$ws.Range("J123").Value = 100000
$ws.Range("L123").Value = 100000
$ws.Range("N123").Value = -109800
$ws.Range("H132").Value = 2954.8865
$ws.Range("I132").Value = 1463.3572
$ws.Range("J132").Value = 34277
$ws.Range("K132").Value = 4390.071599999999
$ws.Range("L132").Value = 102831
$ws.Range("M132").Value = -1860.071599999999
$ws.Range("N132").Value = -107891
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").Value = $null
$ws.Range("N137").Value = 0

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2378.25
$ws.Range("I2").Value = 2446.5
$ws.Range("J2").Value = 2173.5
$ws.Range("K2").Value = 2446.5
$ws.Range("L2").Value = 2173.5
$ws.Range("M2").Value = -2333.5
$ws.Range("N2").Value = -2399.5
$ws.Range("H32").Value = 2747.1372
$ws.Range("I32").Value = 2747.1372
$ws.Range("K32").Value = 2747.1372
$ws.Range("M32").Value = -2460.1372
$ws.Range("H52").Value = 26172.25
$ws.Range("J52").Value = 26172.25
$ws.Range("L52").Value = 26172.25
$ws.Range("N52").Value = -26808.25
$ws.Range("H102").Value = 2104.9473
$ws.Range("I102").Value = 2104.9473
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2104.9473
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = -482.9472999999998
$ws.Range("H110").Value = 1220.5454
$ws.Range("I110").Value = 1220.5454
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1220.5454
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = 824.4546
$ws.Range("H116").Value = 2378.25
$ws.Range("I116").Value = 2446.5
$ws.Range("J116").Value = 2173.5
$ws.Range("K116").Value = 2446.5
$ws.Range("L116").Value = 2173.5
$ws.Range("M116").Value = -152.5
$ws.Range("N116").Value = -6761.5
$ws.Range("H132").Value = 1541.8572
$ws.Range("I132").Value = 1500.4
$ws.Range("J132").Value = 1645.5
$ws.Range("K132").Value = 4501.200000000001
$ws.Range("L132").Value = 4936.5
$ws.Range("M132").Value = -1971.200000000001
$ws.Range("N132").Value = -9996.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2378.25
$ws.Range("I3").Value = 2446.5
$ws.Range("J3").Value = 2173.5
$ws.Range("K3").Value = 2446.5
$ws.Range("L3").Value = 2173.5
$ws.Range("M3").Value = -2332.5
$ws.Range("N3").Value = -2401.5
$ws.Range("H105").Value = 7568.684
$ws.Range("I105").Value = 11920.5
$ws.Range("J105").Value = 2733.3333
$ws.Range("K105").Value = 11920.5
$ws.Range("L105").Value = 2733.3333
$ws.Range("M105").Value = -10173.5
$ws.Range("N105").Value = -6227.3333
$ws.Range("H107").Value = 505000
$ws.Range("I107").Value = 505000
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 505000
$ws.Range("L107").Value = 0
$ws.Range("N107").Value = -503080
$ws.Range("H134").Value = 112469.26
$ws.Range("I134").Value = 150991.1
$ws.Range("J134").Value = 2406.8572
$ws.Range("K134").Value = 452973.3
$ws.Range("L134").Value = 7220.571599999999
$ws.Range("M134").Value = -450438.3
$ws.Range("N134").Value = -12290.5716
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("N137").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4665.3335
$ws.Range("I16").Value = 5098.4
$ws.Range("K16").Value = 5098.4
$ws.Range("M16").Value = -4811.4
$ws.Range("H20").Value = 70000
$ws.Range("I20").Value = 70000
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 70000
$ws.Range("M20").Value = -69764
$ws.Range("N20").Value = 0
$ws.Range("H30").Value = 70000
$ws.Range("I30").Value = 70000
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 70000
$ws.Range("M30").Value = -69909
$ws.Range("N30").Value = 0
$ws.Range("H48").Value = 40050.5
$ws.Range("J48").Value = 40050.5
$ws.Range("L48").Value = 40050.5
$ws.Range("N48").Value = -41002.5
$ws.Range("H99").Value = 4270.7144
$ws.Range("I99").Value = 3979.2
$ws.Range("K99").Value = 3979.2
$ws.Range("M99").Value = -2481.2
$ws.Range("H113").Value = 4665.3335
$ws.Range("I113").Value = 5098.4
$ws.Range("K113").Value = 5098.4
$ws.Range("M113").Value = -2928.4
$ws.Range("H126").Value = 4270.7144
$ws.Range("I126").Value = 3979.2
$ws.Range("K126").Value = 11937.6
$ws.Range("M126").Value = -9467.599999999999
$ws.Range("H128").Value = 70000
$ws.Range("I128").Value = 70000
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 70000
$ws.Range("M128").Value = -65020
$ws.Range("N128").Value = 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 8338.526
$ws.Range("I56").Value = 8338.526
$ws.Range("K56").Value = 8338.526
$ws.Range("M56").Value = -7808.526
$ws.Range("H64").Value = 2399
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = $null
$ws.Range("H67").Value = 2399
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = $null
$ws.Range("H103").Value = 442.14285
$ws.Range("I103").Value = 442.14285
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 1326.42855
$ws.Range("L103").Value = 0
$ws.Range("N103").Value = -447.4285500000001
$ws.Range("H138").Value = 2858.3333
$ws.Range("I138").Value = 1779.25
$ws.Range("K138").Value = 5337.75
$ws.Range("M138").Value = -197.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 23950
$ws.Range("I40").Value = 22900
$ws.Range("J40").Value = 25000
$ws.Range("K40").Value = 22900
$ws.Range("L40").Value = 25000
$ws.Range("M40").Value = -22749
$ws.Range("N40").Value = -25302
$ws.Range("H80").Value = 2703.577
$ws.Range("I80").Value = 1745.2354
$ws.Range("K80").Value = 1745.2354
$ws.Range("M80").Value = -747.2354
$ws.Range("H83").Value = 2703.577
$ws.Range("I83").Value = 1745.2354
$ws.Range("K83").Value = 8726.177
$ws.Range("M83").Value = -3734.177
$ws.Range("H97").Value = 79480.12
$ws.Range("I97").Value = 74519.336
$ws.Range("J97").Value = 85061
$ws.Range("K97").Value = 74519.336
$ws.Range("L97").Value = 85061
$ws.Range("M97").Value = -74023.336
$ws.Range("N97").Value = -86053
$ws.Range("H135").Value = 160880.9
$ws.Range("J135").Value = 170969
$ws.Range("L135").Value = 170969
$ws.Range("N135").Value = -181109
$ws.Range("H140").Value = 114637
$ws.Range("J140").Value = 114637
$ws.Range("L140").Value = 114637
$ws.Range("N140").Value = -124997

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 19998.5
$ws.Range("J42").Value = 19998.5
$ws.Range("L42").Value = 19998.5
$ws.Range("N42").Value = -21124.5
$ws.Range("H49").Value = 19998.5
$ws.Range("J49").Value = 19998.5
$ws.Range("L49").Value = 19998.5
$ws.Range("N49").Value = -20292.5
$ws.Range("H100").Value = 7116.3335
$ws.Range("I100").Value = 7124.5
$ws.Range("K100").Value = 7124.5
$ws.Range("M100").Value = -6583.5
$ws.Range("H122").Value = 115261.11
$ws.Range("I122").Value = 253463
$ws.Range("K122").Value = 760389
$ws.Range("M122").Value = -757939
$ws.Range("H132").Value = 3646.76
$ws.Range("I132").Value = 3217
$ws.Range("K132").Value = 9651
$ws.Range("M132").Value = -7121

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 10135899
$ws.Range("I3").Value = 16759865
$ws.Range("K3").Value = 16759865
$ws.Range("M3").Value = -16759751
$ws.Range("H107").Value = 18519868
$ws.Range("I107").Value = 1496.7222
$ws.Range("K107").Value = 4490.1666
$ws.Range("M107").Value = -2570.1666
